$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "    Creature learning" task row right after "    Creature Names" (row 19) ---
$ws.Rows.Item(20).EntireRow.Insert()
$ws.Range("A20").Value = "    Creature learning"
$ws.Range("AQ20:AS20").Clear()
$ws.Range("AQ19").Copy()
$ws.Range("AT20:AU20").PasteSpecial(-4122)

# --- Insert "        Creature learning" task row right after "        Creature Names" (now row 32) ---
$ws.Rows.Item(33).EntireRow.Insert()
$ws.Range("A33").Value = "        Creature learning"
$ws.Range("AR33").Clear()

# --- Update the remembered selection to match the author's final cursor position ---
$ws.Range("AV33").Select()
